$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 3940
$ws.Range("J10").Value = 3940
$ws.Range("L10").Value = 3940
$ws.Range("N10").Value = -4526
$ws.Range("H40").Value = 2828.5715
$ws.Range("J40").Value = 3360
$ws.Range("L40").Value = 3360
$ws.Range("N40").Value = -3710
$ws.Range("H62").Value = 2901
$ws.Range("I62").Value = 1501.6666
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 1501.6666
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -877.6666
$ws.Range("N62").Value = -6248
$ws.Range("H64").Value = 3925.123
$ws.Range("I64").Value = 3744.1177
$ws.Range("J64").Value = 4123.645
$ws.Range("K64").Value = 3744.1177
$ws.Range("L64").Value = 4123.645
$ws.Range("M64").Value = -3496.1177
$ws.Range("N64").Value = -4619.645
$ws.Range("H65").Value = 2901
$ws.Range("I65").Value = 1501.6666
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 7508.333000000001
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -4388.333000000001
$ws.Range("N65").Value = -31240
$ws.Range("H67").Value = 3925.123
$ws.Range("I67").Value = 3744.1177
$ws.Range("J67").Value = 4123.645
$ws.Range("K67").Value = 3744.1177
$ws.Range("L67").Value = 4123.645
$ws.Range("M67").Value = -2886.1177
$ws.Range("N67").Value = -5839.645
$ws.Range("H132").Value = 3012.3137
$ws.Range("I132").Value = 2517.6
$ws.Range("J132").Value = 4094.5
$ws.Range("K132").Value = 7552.799999999999
$ws.Range("L132").Value = 12283.5
$ws.Range("M132").Value = -5022.799999999999
$ws.Range("N132").Value = -17343.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9108.24
$ws.Range("I32").Value = 6069.6294
$ws.Range("J32").Value = 22062.316
$ws.Range("K32").Value = 6069.6294
$ws.Range("L32").Value = 22062.316
$ws.Range("M32").Value = -5782.6294
$ws.Range("N32").Value = -22636.316
$ws.Range("H74").Value = 1181.9231
$ws.Range("I74").Value = 885.34375
$ws.Range("J74").Value = 2537.7144
$ws.Range("K74").Value = 885.34375
$ws.Range("L74").Value = 2537.7144
$ws.Range("M74").Value = -11.34375
$ws.Range("N74").Value = -4285.7144
$ws.Range("H77").Value = 1181.9231
$ws.Range("I77").Value = 885.34375
$ws.Range("J77").Value = 2537.7144
$ws.Range("K77").Value = 4426.71875
$ws.Range("L77").Value = 12688.572
$ws.Range("M77").Value = -58.71875
$ws.Range("N77").Value = -21424.572
$ws.Range("H88").Value = 1500
$ws.Range("J88").Value = 3000
$ws.Range("L88").Value = 3000
$ws.Range("N88").Value = -3812
$ws.Range("H91").Value = 1500
$ws.Range("J91").Value = 3000
$ws.Range("L91").Value = 3000
$ws.Range("N91").Value = -5808
$ws.Range("H130").Value = 67314.5
$ws.Range("J130").Value = 67314.5
$ws.Range("L130").Value = 67314.5
$ws.Range("N130").Value = -77354.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 14334
$ws.Range("I2").Value = 4502
$ws.Range("J2").Value = 19250
$ws.Range("K2").Value = 4502
$ws.Range("L2").Value = 19250
$ws.Range("M2").Value = -4389
$ws.Range("N2").Value = -19476
$ws.Range("H10").Value = 6405
$ws.Range("I10").Value = 1873.3334
$ws.Range("J10").Value = 20000
$ws.Range("K10").Value = 1873.3334
$ws.Range("L10").Value = 20000
$ws.Range("M10").Value = -1734.3334
$ws.Range("N10").Value = -20278
$ws.Range("H62").Value = 2711.5386
$ws.Range("I62").Value = 2487.5
$ws.Range("J62").Value = 3070
$ws.Range("K62").Value = 2487.5
$ws.Range("L62").Value = 3070
$ws.Range("M62").Value = -1863.5
$ws.Range("N62").Value = -4318
$ws.Range("H65").Value = 2711.5386
$ws.Range("I65").Value = 2487.5
$ws.Range("J65").Value = 3070
$ws.Range("K65").Value = 12437.5
$ws.Range("L65").Value = 15350
$ws.Range("M65").Value = -9317.5
$ws.Range("N65").Value = -21590

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 3200.6667
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 3200.6667
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 9602.000100000001
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -9940.000100000001
$ws.Range("H23").Value = 82.78570999999999
$ws.Range("I23").Value = 110
$ws.Range("J23").Value = 75.36364
$ws.Range("K23").Value = 330
$ws.Range("L23").Value = 226.09092
$ws.Range("M23").Value = -95
$ws.Range("N23").Value = -696.09092
$ws.Range("H27").Value = 3200.6667
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 3200.6667
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 9602.000100000001
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -9806.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 7620.8
$ws.Range("I3").Value = 15150
$ws.Range("J3").Value = 2601.3333
$ws.Range("K3").Value = 15150
$ws.Range("L3").Value = 2601.3333
$ws.Range("M3").Value = -15034
$ws.Range("N3").Value = -2833.3333
$ws.Range("H22").Value = 669.3333
$ws.Range("I22").Value = 1008
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 1008
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -479
$ws.Range("N22").Value = -1558
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H70").Value = 7606.6665
$ws.Range("I70").Value = 7892.857
$ws.Range("J70").Value = 3600
$ws.Range("K70").Value = 7892.857
$ws.Range("L70").Value = 3600
$ws.Range("M70").Value = -7622.857
$ws.Range("N70").Value = -4140
$ws.Range("H73").Value = 7606.6665
$ws.Range("I73").Value = 7892.857
$ws.Range("J73").Value = 3600
$ws.Range("K73").Value = 7892.857
$ws.Range("L73").Value = 3600
$ws.Range("M73").Value = -6956.857
$ws.Range("N73").Value = -5472
$ws.Range("H80").Value = 2884.2444
$ws.Range("I80").Value = 2676.611
$ws.Range("J80").Value = 3022.6667
$ws.Range("K80").Value = 2676.611
$ws.Range("L80").Value = 3022.6667
$ws.Range("M80").Value = -1678.611
$ws.Range("N80").Value = -5018.6667
$ws.Range("H83").Value = 2884.2444
$ws.Range("I83").Value = 2676.611
$ws.Range("J83").Value = 3022.6667
$ws.Range("K83").Value = 13383.055
$ws.Range("L83").Value = 15113.3335
$ws.Range("M83").Value = -8391.055
$ws.Range("N83").Value = -25097.3335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 2150.5
$ws.Range("I4").Value = 500
$ws.Range("J4").Value = 2480.6
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = 2480.6
$ws.Range("M4").Value = -387
$ws.Range("N4").Value = -2706.6
